{"js": "// Apply the \"Game Overview\" content edits described in the commit.\n//\n// Strategy: locate each changed phrase with Body.search() (exact,\n// case-sensitive literal text) and rewrite it in place with\n// Range.insertText(..., \"Replace\"). This preserves every paragraph\n// that is untouched and only rewrites the text that actually changed.\n//\n// The one structural change \u2014 the \"_GoBack\" bookmark moving from the\n// \"Game Objective\" paragraph (where it used to sit between \"The ob\"\n// and \"jective...\") into the title paragraph (between \"Lost\" and\n// \" Dreams\") \u2014 is replicated explicitly at the end.\n\nconst body = context.document.body;\n\n// Find the unique occurrence of `needle` and replace it with\n// `replacement`. Throws if the match count isn't exactly 1 so a typo\n// in the needle text surfaces loudly instead of silently editing the\n// wrong spot (or every spot).\nasync function replaceOnce(needle, replacement) {\n  const results = body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(needle) + \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) Title: \"Game Overview: Finding Dreams\" -> \"Game Overview: Lost Dreams\"\nawait replaceOnce(\"Finding\", \"Lost\");\n\n// 2) Game Story paragraph tweaks.\nawait replaceOnce(\"The game portray a\", \"The game portrays a\");\nawait replaceOnce(\"a hard days work\", \"a hard day\\u2019s work\");\nawait replaceOnce(\n  \"counting his sheep knowing that one his sheep\",\n  \"counting his sheep knowing that once his sheep\"\n);\n\n// 3) Gameplay Mechanics / 2D paragraph tweaks.\nawait replaceOnce(\"but with a key difference. Instead\", \"but with a key difference: Instead\");\nawait replaceOnce(\"flat or bumpy terrain, the train is made up\", \"flat or bumpy terrain, the terrain is made up\");\nawait replaceOnce(\"move around the parimiter of the shape of the island\", \"move around the perimeter of the island\");\nawait replaceOnce(\"This not only introduces  a novel visual\", \"This not only introduces a novel visual\");\nawait replaceOnce(\"jump between planets to travers a level\", \"jump between islands to traverse a level\");\nawait replaceOnce(\n  \"each level would consist of multiple, different, islands.\",\n  \"each level will consist of multiple unique islands.\"\n);\n\n// 4) Player / enemies paragraph tweaks.\nawait replaceOnce(\n  \"in a level without touching any enemies. Touching an enemy means that the player has failed and must restart the level (the level resetting to its original state).\",\n  \"in a level while avoiding enemies. When the player makes contact with an enemy, one of the player\\u2019s sheep is lost and put back somewhere in the world for the player to retrieve again.\"\n);\nawait replaceOnce(\n  \"Each star collected and sheep collected (subject to multipliers) adds points.\",\n  \"Each star and sheep collected adds points (subject to multipliers).\"\n);\n\n// 5) Move the \"_GoBack\" bookmark out of the \"Game Objective\" paragraph\n// and into the title, now landing right after \"Lost\".\n//\n// Document.deleteBookmark removes just the bookmark markers (unlike\n// Range.delete() on the bookmark's \u2014 empty \u2014 range, which would blow\n// away the whole remaining paragraph text), so use that instead of\n// deleting body.getBookmarkRange(\"_GoBack\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst lostResults = body.search(\"Lost\", { matchCase: true });\nlostResults.load(\"text\");\nawait context.sync();\nif (lostResults.items.length !== 1) {\n  throw new Error(\"Expected exactly 1 match for \\\"Lost\\\" but found \" + lostResults.items.length);\n}\nconst afterLost = lostResults.items[0].getRange(\"End\");\nafterLost.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Apply the \"Game Overview\" content edits described in the commit.\n#\n# Strategy: use Range.Find/Replacement (Word's Find & Replace engine) to\n# rewrite each changed phrase in place. Every call targets a long,\n# unique literal substring so exactly one (correct) location is\n# rewritten per call.\n#\n# The one structural change \u2014 the \"_GoBack\" bookmark moving from the\n# \"Game Objective\" paragraph (where it used to sit between \"The ob\"\n# and \"jective...\") into the title paragraph (between \"Lost\" and\n# \" Dreams\") \u2014 is replicated explicitly at the end.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param($FindText, $ReplaceText)\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n    # wdReplaceAll = 2; each FindText above is crafted to be unique in\n    # the document so this only ever touches the one intended spot.\n    $ok = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    if (-not $ok) {\n        throw \"Replace failed for: $FindText\"\n    }\n}\n\n# 1) Title: \"Game Overview: Finding Dreams\" -> \"Game Overview: Lost Dreams\"\nReplace-Text \"Finding\" \"Lost\"\n\n# 2) Game Story paragraph tweaks.\nReplace-Text \"The game portray a\" \"The game portrays a\"\nReplace-Text \"a hard days work\" \"a hard day\u2019s work\"\nReplace-Text \"counting his sheep knowing that one his sheep\" \"counting his sheep knowing that once his sheep\"\n\n# 3) Gameplay Mechanics / 2D paragraph tweaks.\nReplace-Text \"but with a key difference. Instead\" \"but with a key difference: Instead\"\nReplace-Text \"flat or bumpy terrain, the train is made up\" \"flat or bumpy terrain, the terrain is made up\"\nReplace-Text \"move around the parimiter of the shape of the island\" \"move around the perimeter of the island\"\nReplace-Text \"This not only introduces  a novel visual\" \"This not only introduces a novel visual\"\nReplace-Text \"jump between planets to travers a level\" \"jump between islands to traverse a level\"\nReplace-Text \"each level would consist of multiple, different, islands.\" \"each level will consist of multiple unique islands.\"\n\n# 4) Player / enemies paragraph tweaks.\nReplace-Text \"in a level without touching any enemies. Touching an enemy means that the player has failed and must restart the level (the level resetting to its original state).\" \"in a level while avoiding enemies. When the player makes contact with an enemy, one of the player\u2019s sheep is lost and put back somewhere in the world for the player to retrieve again.\"\nReplace-Text \"Each star collected and sheep collected (subject to multipliers) adds points.\" \"Each star and sheep collected adds points (subject to multipliers).\"\n\n# 5) Move the \"_GoBack\" bookmark out of the \"Game Objective\" paragraph\n# and into the title, now landing right after \"Lost\".\n#\n# Bookmark.Delete() removes just the bookmark markers without touching\n# the surrounding run text.\n$oldBookmark = $d.Bookmarks.Item(\"_GoBack\")\n$oldBookmark.Delete()\n\n$titleFind = $d.Content.Find\n$titleFind.ClearFormatting()\n$titleFind.Text = \"Lost\"\n$titleFind.MatchCase = $true\n$titleFind.MatchWildcards = $false\n$titleFind.Forward = $true\n$titleFind.Wrap = 1\n$titleFind.Execute() | Out-Null\nif (-not $titleFind.Found) {\n    throw \"Could not find 'Lost' in the title to re-anchor the _GoBack bookmark\"\n}\n\n$lostRange = $d.Content\n$lostRange.Find.Execute(\"Lost\") | Out-Null\n$lostRange.Collapse(0)  # wdCollapseEnd = 0\n$d.Bookmarks.Add(\"_GoBack\", $lostRange)\n\nWrite-Output \"Game Overview edits applied\"\n"}
